$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 496
$ws1.Range("F6").Value = 14321
$ws1.Range("F7").Value = 16440
$ws1.Range("F22").Value = 134
$ws1.Range("F25").Value = 2
$ws1.Range("F32").Value = 5733

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 496
$ws4.Range("F6").Value = 14321
$ws4.Range("F7").Value = 16440
$ws4.Range("F22").Value = 134
$ws4.Range("F26").Value = 2
$ws4.Range("F35").Value = 5733
